$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "TEMPERATURA" is captured by Selenium as plain text (e.g. "28"), so the
# cell must keep it as text rather than letting Excel auto-convert it to a
# number. Pre-format that single cell as Text before writing the value.
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "Wednesday 10:00"
$ws.Range("B2").Value = "28"
$ws.Range("C2").Value = "Humidity: 46%"
